# Weekly price-report update: a new daily record is inserted at the top of
# the data block (row 221), pushing every existing record down by one row;
# the record that used to be last (old row 321) ends up at the new row 322.
#
# New data point added (row 221 after the insert):
#   Fecha (D) = 44636, Volumen (J) = 140, Precio minimo (K) = 4500,
#   Precio maximo (L) = 4800, Precio promedio ponderado (M) = 4629,
#   Precio $/Kg (P) = 1543
# All the other columns of the new record repeat the values that were
# already constant across the whole "Cilantro" block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 221; this shifts rows 221:321 down to 222:322
# and grows the sheet dimension to A1:R322 automatically.
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with the new record.
$ws.Range("A221").Value = 3
$ws.Range("B221").Value = "Femacal de La Calera"
$ws.Range("C221").Value = "Coquimbo"
$ws.Range("D221").Value = 44636
$ws.Range("E221").Value = 5
$ws.Range("F221").Value = 100112040
$ws.Range("G221").Value = "Cilantro"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 140
$ws.Range("K221").Value = 4500
$ws.Range("L221").Value = 4800
$ws.Range("M221").Value = 4629
$ws.Range("N221").Value = "$/docena de atados (3 kilos)"
$ws.Range("O221").Value = "Provincia de Quillota"
$ws.Range("P221").Value = 1543
$ws.Range("Q221").Value = 3
$ws.Range("R221").Value = "Hortaliza"

# Preserve the same numeric format the other "Fecha" cells use (date-style
# display for the serial date value), matching style index already carried
# over to the inserted row by Insert().
$ws.Range("D221").NumberFormat = $ws.Range("D222").NumberFormat
